$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","L1cam","Egfr","ECs",3,1,11.88712433333333,35.661373,0.5967229292030898,0.5967229292030898,3,1,1.854571666666667,5.563715,0.01651371646154392,0.01651371646154392,22.04552398674389,198.409715880695,0.00985411325896177,0.009854113258961772),
  @("ECs","L1cam","Egfr","FAPs",3,1,11.88712433333333,35.661373,0.5967229292030898,0.5967229292030898,3,1,92.91372433333333,278.741173,0.8273343794712995,0.8273343794712996,1104.476993423392,9940.292940810528,0.4936893943485345,0.4936893943485345),
  @("ECs","L1cam","Egfr","M2",3,1,11.88712433333333,35.661373,0.5967229292030898,0.5967229292030898,3,1,0.162136,0.4864080000000001,0.001443712303133186,0.001443712303133187,1.927330790909334,17.345977118184,0.000861496234452174,0.0008614962344521742),
  @("ECs","L1cam","Egfr","sCs",3,1,11.88712433333333,35.661373,0.5967229292030898,0.5967229292030898,3,1,17.37449166666667,52.123475,0.1547081917640233,0.1547081917640233,206.5327426701305,1858.794684031175,0.09231792536114132,0.09231792536114132),
  @("FAPs","L1cam","Egfr","ECs",2,0.6666666666666666,0.3785896666666667,1.135769,0.01900486009268527,0.01900486009268527,3,1,1.854571666666667,5.563715,0.01651371646154392,0.01651371646154392,0.7021216690927777,6.319095021835,0.0003138408709619159,0.000313840870961916),
  @("FAPs","L1cam","Egfr","FAPs",2,0.6666666666666666,0.3785896666666667,1.135769,0.01900486009268527,0.01900486009268527,3,1,92.91372433333333,278.741173,0.8273343794712995,0.8273343794712996,35.17617592411522,316.585583317037,0.01572337413172064,0.01572337413172064),
  @("FAPs","L1cam","Egfr","M2",2,0.6666666666666666,0.3785896666666667,1.135769,0.01900486009268527,0.01900486009268527,3,1,0.162136,0.4864080000000001,0.001443712303133186,0.001443712303133187,0.06138301419466668,0.5524471277520001,0.00002743755033513464,0.00002743755033513464),
  @("FAPs","L1cam","Egfr","sCs",2,0.6666666666666666,0.3785896666666667,1.135769,0.01900486009268527,0.01900486009268527,3,1,17.37449166666667,52.123475,0.1547081917640233,0.1547081917640233,6.577803008586112,59.200227077275,0.002940207539667587,0.002940207539667587),
  @("M2","L1cam","Egfr","ECs",3,1,4.123197999999999,12.369594,0.2069808238940481,0.2069808238940482,3,1,1.854571666666667,5.563715,0.01651371646154392,0.01651371646154392,7.646766186856666,68.82089568171,0.003418022638763066,0.003418022638763067),
  @("M2","L1cam","Egfr","FAPs",3,1,4.123197999999999,12.369594,0.2069808238940481,0.2069808238940482,3,1,92.91372433333333,278.741173,0.8273343794712995,0.8273343794712996,383.1016823437513,3447.915141093762,0.1712423514988407,0.1712423514988407),
  @("M2","L1cam","Egfr","M2",3,1,4.123197999999999,12.369594,0.2069808238940481,0.2069808238940482,3,1,0.162136,0.4864080000000001,0.001443712303133186,0.001443712303133187,0.668518830928,6.016669478352,0.0002988207619684807,0.0002988207619684808),
  @("M2","L1cam","Egfr","sCs",3,1,4.123197999999999,12.369594,0.2069808238940481,0.2069808238940482,3,1,17.37449166666667,52.123475,0.1547081917640233,0.1547081917640233,71.63846929101666,644.74622361915,0.03202162899447594,0.03202162899447594),
  @("sCs","L1cam","Egfr","ECs",3,1,3.531764333333333,10.595293,0.1772913868101768,0.1772913868101768,3,1,1.854571666666667,5.563715,0.01651371646154392,0.01651371646154392,6.549910065943889,58.949190593495,0.002927739692857166,0.002927739692857167),
  @("sCs","L1cam","Egfr","FAPs",3,1,3.531764333333333,10.595293,0.1772913868101768,0.1772913868101768,3,1,92.91372433333333,278.741173,0.8273343794712995,0.8273343794712996,328.1493776776321,2953.344399098689,0.1466792594922037,0.1466792594922037),
  @("sCs","L1cam","Egfr","M2",3,1,3.531764333333333,10.595293,0.1772913868101768,0.1772913868101768,3,1,0.162136,0.4864080000000001,0.001443712303133186,0.001443712303133187,0.5726261419493335,5.153635277544001,0.0002559577563773969,0.0002559577563773969),
  @("sCs","L1cam","Egfr","sCs",3,1,3.531764333333333,10.595293,0.1772913868101768,0.1772913868101768,3,1,17.37449166666667,52.123475,0.1547081917640233,0.1547081917640233,61.36260997813056,552.2634898031749,0.02742842986873846,0.02742842986873846)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $excelRow = $i + 2
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
  }
}
